$d = $word.ActiveDocument

$para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "Data Engineering and Infrastructure Architecture") {
        $para = $p
        break
    }
}

$r1 = $para.Range
$r1.Collapse(0)
$r1.InsertParagraphAfter()
$r1.InsertAfter("• Architected data infrastructure processing 15+ billion voter records to support meta-analytical voter file corrections")

$para2 = $para.Next()
$r2 = $para2.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$r2.InsertAfter("• Built scalable ETL pipelines enabling analysis of 50,000+ electoral boundaries across all levels of government")

$para3 = $para2.Next()
$r3 = $para3.Range
$r3.Collapse(0)
$r3.InsertParagraphAfter()
$r3.InsertAfter("• Developed Python boundary estimation algorithm that reduced mapping costs by 75% for 200+ organizations")
